$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Jennifer%Lighter%Jennifer.Lighter@nyumc.org%1,                  Michael%Phillips%NULL%1,                  Sarah%Hochman%NULL%1,                  Stephanie%Sterling%NULL%1,                  Diane%Johnson%NULL%1,                  Fritz%Francois%NULL%0,                  Anna%Stachel%NULL%1]"
$ws.Range("I2").Value = "_PMC_CROSSREF"
$ws.Range("E3").Value = "[J.%Wu%xref no email%1,            W.%Li%xref no email%1,            X.%Shi%xref no email%1,            Z.%Chen%xref no email%1,            B.%Jiang%xref no email%1,            J.%Liu%xref no email%1,            D.%Wang%xref no email%1,            C.%Liu%xref no email%1,            Y.%Meng%xref no email%1,            L.%Cui%xref no email%1,            J.%Yu%xref no email%1,            H.%Cao%xref no email%1,            L.%Li%xref no email%1]"
$ws.Range("I3").Value = "_CROSSREF"
